$wb = $excel.ActiveWorkbook

# --- Sheet references ---
$wsMeet = $wb.Worksheets.Item("MeetResultaat")
$wsTHP = $wb.Worksheets.Item("THP")
$wsPrem = $wb.Worksheets.Item("isPrematuur")

# --- MeetResultaat: rename C1/D1/E1 headers first ---
$wsMeet.Range("C1").Value = "HPS.T4"
$wsMeet.Range("D1").Value = "HPS.TBG"
$wsMeet.Range("E1").Value = "HPS.TSH"

# --- THP: update the field-name label in A1 ---
$wsTHP.Range("A1").Value = "HPS.hielprik.hielprikType"

# --- MeetResultaat: rename the "hps.redenOnbekend(...)" outcome labels used throughout
#     column G to the new "HPS.bruikbaarheid(...)" naming ---
$rowsTSH = @(5,15,16,17,37,39,40,41,42,43,45,54,62,76)
foreach ($r in $rowsTSH) {
    $wsMeet.Cells.Item($r, 7).Value = 'HPS.bruikbaarheid(HPS.hielprik,"TSH")'
}

$rowsT4TBG = @(29,30,44)
foreach ($r in $rowsT4TBG) {
    $wsMeet.Cells.Item($r, 7).Value = 'HPS.bruikbaarheid(HPS.hielprik,"T4/TBG")'
}

$rowsT4 = @(77,78,79,81,82)
foreach ($r in $rowsT4) {
    $wsMeet.Cells.Item($r, 7).Value = 'HPS.bruikbaarheid(HPS.hielprik,"T4")'
}

# --- MeetResultaat: rename B1 (isTweedeHielprik) and F1 (HPS.T4_TBG), and A1 (isPrematuur) ---
$wsMeet.Range("A1").Value = "isPrematuur"
$wsMeet.Range("B1").Value = "isTweedeHielprik"
$wsMeet.Range("F1").Value = "HPS.T4_TBG"

# --- Selection / active sheet state: MeetResultaat becomes the active tab, cell F2 selected ---
$wsMeet.Activate()
$wsMeet.Range("F2").Select()
